$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-29 Wednesday", "2024-05-30 Thursday"),
    @("51×90=", "79×60="),
    @("66×24=", "22×31="),
    @("89×33=", "94×31="),
    @("25×48=", "26×42="),
    @("16×97=", "17×72="),
    @("11×58=", "51×53="),
    @("47×91=", "41×61="),
    @("65×95=", "51×65="),
    @("70×96=", "33×50="),
    @("76×62=", "96×36="),
    @("49×99=", "86×29="),
    @("65×15=", "13×73="),
    @("62×30=", "48×56="),
    @("95×97=", "53×64="),
    @("22×85=", "17×71="),
    @("23×81=", "99×73="),
    @("60×88=", "75×41="),
    @("15×62=", "29×26="),
    @("15×26=", "28×63="),
    @("26×46=", "49×81="),
    @("33×66=", "36×43="),
    @("96×44=", "78×79="),
    @("79×27=", "66×91="),
    @("60×75=", "12×77="),
    @("25×69=", "12×27=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
